$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column AD (30): everything from AD onward
# (the "nom" / "url_produit" columns) shifts one column to the right,
# becoming AE / AF. This makes room for a new timestamped price-history
# column at AD, mirroring the existing per-run price columns B..AC.
$ws.Columns.Item(30).Insert()

# Header for the newly inserted column (row 1) — same style as the other
# timestamp header cells in row 1.
$ws.Range("AD1").Value = "2026-01-28 23:16:49"

# Populate the new AD column (rows 2-80) with the latest price snapshot,
# copied from column AC (the last existing price-history column) for
# those same rows. Rows 81-205 have no price data (AC is blank there),
# so AD stays blank for them too.
$ws.Range("AD2:AD80").Value2 = $ws.Range("AC2:AC80").Value2
